$d = $word.ActiveDocument
$r = $d.Content
$found = $r.Find.Execute("its enemy targeted on its board, and saves", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output $found
Write-Output $d.Paragraphs(5).Range.Text
